$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (this also updates the _FilterDatabase defined name reference)
$ws.Name = "Jul 2020 to Aug 2020"

# Move the selection to A23
$ws.Range("A23").Select()
